# FunctionDB.xlsx - extend the B/C/D "function id" table down to row 61
# and fill in the missing D values for rows 29-34 (columns: strcn / cn / dr).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill the few D-column gaps (rows 29-34) so the D series matches B/C ---
$ws.Range("D29").Value = 28
$ws.Range("D30").Value = 29
$ws.Range("D31").Value = 30
$ws.Range("D32").Value = 31
$ws.Range("D33").Value = 32
$ws.Range("D34").Value = 33

# --- Append rows 47-61, continuing the B/C numbering sequence (46-60) ---
$newRows = @(
    @(47, 46), @(48, 47), @(49, 48), @(50, 49), @(51, 50),
    @(52, 51), @(53, 52), @(54, 53), @(55, 54), @(56, 55),
    @(57, 56), @(58, 57), @(59, 58), @(60, 59), @(61, 60)
)

foreach ($pair in $newRows) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Cells.Item($r, 2).Value = $v
    $ws.Cells.Item($r, 3).Value = $v
}

# --- Match the final selection left by the author (B59:C61, active B59) ---
$null = $ws.Range("B59:C61").Select()
